$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Batch")

# Insert two new rows before row 21, copying the format/style of the
# "Scheme"/"Subsidy" rows (19:20) so that the new rows inherit the same
# look (font, row style) as their neighbours.
$ws.Rows("19:20").Copy()
$ws.Rows("21:22").Insert(-4121)
$ws.Application.CutCopyMode = $false

# Row 17/18 used to read "Vaccination Date" / "Vaccination Report
# Generation Date"; rename them to the "[1st] ..." variants and reuse
# the old text for the newly created "[2nd ...]" rows below the
# Scheme/Subsidy pair.
$ws.Range("A17").Value = "[1st] Vaccination Date"
$ws.Range("A18").Value = "[1st] Vaccination Report Generation Date"

$ws.Range("A19").Value = "[2nd Vaccination Date]"
$ws.Range("A20").Value = "[2nd Vaccination Report Generation Date]"

$ws.Range("A21").Value = "Scheme"
$ws.Range("A22").Value = "Subsidy"
